# chore: update Sheets via scheduled runner
# Refresh cached market-price figures (columns H-N) for the affected leve rows
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR profit sheets.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 371676.12
$ws.Range("J129").Value = 477766.44
$ws.Range("L129").Value = 1433299.32
$ws.Range("N129").Value = -1443299.32
$ws.Range("H132").Value = 3026.1936
$ws.Range("I132").Value = 3546.5417
$ws.Range("K132").Value = 10639.6251
$ws.Range("M132").Value = -8109.625100000001
$ws.Range("H137").Value = 1743.4054
$ws.Range("J137").Value = 1771.8572
$ws.Range("L137").Value = 5315.571599999999
$ws.Range("N137").Value = -10415.5716
$ws.Range("H141").Value = 2097.8
$ws.Range("I141").Value = 1211.4286
$ws.Range("J141").Value = 4166
$ws.Range("K141").Value = 3634.2858
$ws.Range("L141").Value = 12498
$ws.Range("M141").Value = 1545.7142
$ws.Range("N141").Value = -22858

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3175.4827
$ws.Range("I32").Value = 2793.1013
$ws.Range("J32").Value = 6951.5
$ws.Range("K32").Value = 2793.1013
$ws.Range("L32").Value = 6951.5
$ws.Range("M32").Value = -2506.1013
$ws.Range("N32").Value = -7525.5
$ws.Range("H61").Value = 451411.78
$ws.Range("I61").Value = 487553.3
$ws.Range("J61").Value = 5666.3335
$ws.Range("K61").Value = 487553.3
$ws.Range("L61").Value = 5666.3335
$ws.Range("M61").Value = -487341.3
$ws.Range("N61").Value = -6090.3335
$ws.Range("H74").Value = 40002476
$ws.Range("I74").Value = 45457040
$ws.Range("J74").Value = 2333
$ws.Range("K74").Value = 45457040
$ws.Range("L74").Value = 2333
$ws.Range("M74").Value = -45456166
$ws.Range("N74").Value = -4081
$ws.Range("H77").Value = 40002476
$ws.Range("I77").Value = 45457040
$ws.Range("J77").Value = 2333
$ws.Range("K77").Value = 227285200
$ws.Range("L77").Value = 11665
$ws.Range("M77").Value = -227280832
$ws.Range("N77").Value = -20401
$ws.Range("H136").Value = 451411.78
$ws.Range("I136").Value = 487553.3
$ws.Range("J136").Value = 5666.3335
$ws.Range("K136").Value = 1462659.9
$ws.Range("L136").Value = 16999.0005
$ws.Range("M136").Value = -1460109.9
$ws.Range("N136").Value = -22099.0005

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2305.647
$ws.Range("I134").Value = 2338.5305
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 7015.5915
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -4480.5915
$ws.Range("N134").Value = -9570

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3742.578
$ws.Range("I31").Value = 2692.8333
$ws.Range("K31").Value = 2692.8333
$ws.Range("M31").Value = -2397.8333
$ws.Range("H34").Value = 3742.578
$ws.Range("I34").Value = 2692.8333
$ws.Range("K34").Value = 2692.8333
$ws.Range("M34").Value = -2490.8333
$ws.Range("H58").Value = 12792.116
$ws.Range("I58").Value = 965.7742
$ws.Range("K58").Value = 965.7742
$ws.Range("M58").Value = -762.7742
$ws.Range("H99").Value = 15154174
$ws.Range("I99").Value = 2341.0476
$ws.Range("K99").Value = 2341.0476
$ws.Range("M99").Value = -843.0475999999999
$ws.Range("H126").Value = 15154174
$ws.Range("I126").Value = 2341.0476
$ws.Range("K126").Value = 7023.1428
$ws.Range("M126").Value = -4553.1428
$ws.Range("H132").Value = 2276.4146
$ws.Range("I132").Value = 1587.3948
$ws.Range("K132").Value = 4762.1844
$ws.Range("M132").Value = -2232.1844
$ws.Range("H134").Value = 865.75
$ws.Range("I134").Value = 722.23334
$ws.Range("J134").Value = 1583.3334
$ws.Range("K134").Value = 2166.70002
$ws.Range("L134").Value = 4750.0002
$ws.Range("M134").Value = 368.2999799999998
$ws.Range("N134").Value = -9820.0002
$ws.Range("H136").Value = 12792.116
$ws.Range("I136").Value = 965.7742
$ws.Range("K136").Value = 2897.3226
$ws.Range("M136").Value = -347.3226

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 136
$ws.Range("I50").Value = 114.666664
$ws.Range("K50").Value = 343.999992
$ws.Range("M50").Value = 137.000008
$ws.Range("H53").Value = 136
$ws.Range("I53").Value = 114.666664
$ws.Range("K53").Value = 343.999992
$ws.Range("M53").Value = 137.000008
$ws.Range("H74").Value = 10011.25
$ws.Range("J74").Value = 10011.25
$ws.Range("L74").Value = 30033.75
$ws.Range("N74").Value = -32155.75
$ws.Range("H77").Value = 10011.25
$ws.Range("J77").Value = 10011.25
$ws.Range("L77").Value = 90101.25
$ws.Range("N77").Value = -100709.25
$ws.Range("H116").Value = 1768.2222
$ws.Range("J116").Value = 1737.7142
$ws.Range("L116").Value = 5213.142599999999
$ws.Range("N116").Value = -12097.1426
$ws.Range("H123").Value = 3200.6191
$ws.Range("I123").Value = 2481.75
$ws.Range("J123").Value = 3369.7646
$ws.Range("K123").Value = 7445.25
$ws.Range("L123").Value = 10109.2938
$ws.Range("M123").Value = -4995.25
$ws.Range("N123").Value = -15009.2938
$ws.Range("H129").Value = 1410.6296
$ws.Range("I129").Value = 961.8182
$ws.Range("J129").Value = 1719.1875
$ws.Range("K129").Value = 2885.4546
$ws.Range("L129").Value = 5157.5625
$ws.Range("M129").Value = 2114.5454
$ws.Range("N129").Value = -15157.5625
$ws.Range("H130").Value = 2999
$ws.Range("J130").Value = 2999
$ws.Range("L130").Value = 8997
$ws.Range("N130").Value = -19037
$ws.Range("H131").Value = 294936.66
$ws.Range("J131").Value = 371225.4
$ws.Range("L131").Value = 1113676.2
$ws.Range("N131").Value = -1123756.2
$ws.Range("H134").Value = 2519.1904
$ws.Range("I134").Value = 1863.75
$ws.Range("J134").Value = 4616.6
$ws.Range("K134").Value = 5591.25
$ws.Range("L134").Value = 13849.8
$ws.Range("M134").Value = -521.25
$ws.Range("N134").Value = -23989.8
$ws.Range("H136").Value = 2193.1
$ws.Range("J136").Value = 4283.6665
$ws.Range("L136").Value = 12850.9995
$ws.Range("N136").Value = -23050.9995
$ws.Range("H137").Value = 30306938
$ws.Range("I137").Value = 1030
$ws.Range("J137").Value = 33337530
$ws.Range("K137").Value = 3090
$ws.Range("L137").Value = 100012590
$ws.Range("M137").Value = 2010
$ws.Range("N137").Value = -100022790
$ws.Range("H138").Value = 1361.3572
$ws.Range("I138").Value = 1361.3572
$ws.Range("K138").Value = 4084.0716
$ws.Range("M138").Value = 1055.9284
$ws.Range("H139").Value = 1803.8889
$ws.Range("I139").Value = 1203.5454
$ws.Range("K139").Value = 3610.6362
$ws.Range("M139").Value = 1529.3638
$ws.Range("H140").Value = 3525.4167
$ws.Range("I140").Value = 1832.8572
$ws.Range("J140").Value = 5895
$ws.Range("K140").Value = 5498.571599999999
$ws.Range("L140").Value = 17685
$ws.Range("M140").Value = -318.5715999999993
$ws.Range("N140").Value = -28045
$ws.Range("H141").Value = 3611
$ws.Range("I141").Value = 3611
$ws.Range("K141").Value = 10833
$ws.Range("M141").Value = -5653

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3922.6924
$ws.Range("I132").Value = 2374.625
$ws.Range("K132").Value = 7123.875
$ws.Range("M132").Value = -4593.875
$ws.Range("H136").Value = 863.1515000000001
$ws.Range("I136").Value = 863.1515000000001
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2589.4545
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -39.45450000000028
$ws.Range("N136").ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1114.4286
$ws.Range("I132").Value = 740.2727
$ws.Range("J132").Value = 2486.3333
$ws.Range("K132").Value = 2220.8181
$ws.Range("L132").Value = 7458.999899999999
$ws.Range("M132").Value = 309.1819
$ws.Range("N132").Value = -12518.9999

Write-Output "applied price refresh to ALC, ARM, BSM, CRP, CUL, LTW, WVR"
